$d = $word.ActiveDocument

# 1) Merge the two runs describing GDB command definitions into one,
#    fixing the word break ("definitio" + "ns):" -> "definitions):").
$d.Content.Find.Execute(
    "Describe what each of the following GDB commands does (experiment from within GDB or use the internet to find their definitio" + "ns):",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Describe what each of the following GDB commands does (experiment from within GDB or use the internet to find their definitions):",
    2)

# 2) Replace the "Demonstrate the working program..." sentence: drop the
#    mailto hyperlink to Cole.Scott.Peterson@huskers.unl.edu and change the
#    closing instructions to point to "handin" instead.
$d.Content.Find.Execute(
    "please send this completed worksheet to Cole.Scott.Peterson@huskers.unl.edu to verify completion.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "please submit this completed worksheet to handin.",
    2)

# The newly-inserted closing phrase gets a slightly larger font size (12pt)
# as its own run.
$rng = $d.Content
$rng.Find.Execute("please submit this completed worksheet to handin.")
$rng.Font.Size = 12
